$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Orçamento do projeo"

$ws.Range("A3").Value = "Item de gasto"
$ws.Range("A5").Value = "Salários"
$ws.Range("A6").Value = "Telefonia"

$ws.Range("B3").Value = "Valor mensal previsto"
$ws.Range("B5").Value = "R$ 120.000"
$ws.Range("B6").Value = "R$ 34.000"

$ws.Columns.Item(1).ColumnWidth = 19.3
$ws.Columns.Item(2).ColumnWidth = 19.6

$ws.Range("B10").Select() | Out-Null
